# "Some database tests added."
# Adds a new data row (Code=4, Kind="Company") below the existing table,
# then leaves the selection on the next empty cell (A6), mirroring the
# cursor position Excel leaves after typing the new row and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Company"

$ws.Range("A6").Select()
